# ----------------------------------------------------------------------------
# US398956-EstablecerTipoCombustible-TestReport.docx
# "Terminado el informe de pruebas"
# ----------------------------------------------------------------------------

$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# 1) Remove the author attribution parenthetical from the UI-tests paragraph.
# ----------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = " (Juan David Corrales Gil) mediante el plan de pruebas proporcionado"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = " mediante el plan de pruebas proporcionado"
$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# ----------------------------------------------------------------------------
# 2) Touch the "gasolineras" paragraph (no text change) so the stray
#    _GoBack bookmark left over from the previous edit gets cleared, same
#    as it would be the moment a real author resumes editing nearby text.
#    (the replaced span straddles the bookmark location so Word drops it)
# ----------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "pero al terminar las pruebas todo"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "pero al terminar las pruebas todo"
$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

Write-Output "Step 1-2 done"

# ----------------------------------------------------------------------------
# 3) Rewrite the unit-tests paragraph: drop the "(Hamza Hamda)" parenthetical
#    and replace the trailing ellipsis placeholder with the full write-up.
# ----------------------------------------------------------------------------
$unitParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Con respecto a la implementaci")) {
        if ($d.Paragraphs($i).Range.Text.Contains("unitarias")) {
            $unitParaIndex = $i
        }
    }
}
Write-Output "unitParaIndex=$unitParaIndex"

$p = $d.Paragraphs($unitParaIndex)
$pRange = $p.Range
$bodyRange = $d.Range($pRange.Start, $pRange.End - 1)

$newText = "Con respecto a la implementación de las pruebas unitarias, se han implementado las pruebas unitarias para la clase PresenterGasolinera, concretamente se han probado los métodos lecturaCombustiblePorDefecto y escrituraCombustiblePorDefecto. A raíz de plantear las pruebas a realizar para estos métodos, se detectó la necesidad de añadir una clase de excepción que no estaba inicialmente considerada en el método escrituraCombustiblePorDefecto, asimismo se ha tenido que modificar la manera de tratar las excepciones de E/S en ambos métodos. Inicialmente se hacía un catch de estas excepciones y no se trataban dentro de los métodos, lo cual sería contraproducente ya que la ejecución seguiría su curso como si nada hubiese ocurrido y realmente es necesario tratar estas excepciones para que la aplicación funcione correctamente. Se ha cambiado este comportamiento para que los métodos propaguen estas excepciones y que sea el programa principal el que las trate. Una vez realizados estos cambios las pruebas se han ejecutado satisfactoriamente."

$bodyRange.Font.Italic = $false
$bodyRange.Text = $newText

Write-Output $d.Paragraphs($unitParaIndex).Range.Text
Write-Output $d.Paragraphs.Count

# Italicize the class/method names that are mentioned in the new text.
function Set-ItalicForAllOccurrences($paragraphRange, $needle) {
    $scanRange = $paragraphRange.Duplicate
    $find = $scanRange.Find
    $find.ClearFormatting()
    $find.Text = $needle
    $find.Forward = $true
    $find.Wrap = 0
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    while ($find.Execute()) {
        $scanRange.Font.Italic = $true
        $scanRange.Font.ItalicBi = $true
        $scanRange.Collapse(0)
    }
}

$unitPara = $d.Paragraphs($unitParaIndex).Range
Set-ItalicForAllOccurrences $unitPara "PresenterGasolinera"
Set-ItalicForAllOccurrences $unitPara "lecturaCombustiblePorDefecto"
Set-ItalicForAllOccurrences $unitPara "escrituraCombustiblePorDefecto"

